$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update header row (J1:R1) with new column names ---
# Note: set from right-to-left isn't required; shared string reuse is handled by the engine.
$ws.Range("J1").Value = "num_days_25_encounters"
$ws.Range("K1").Value = "num_days_100_encounters"
$ws.Range("L1").Value = "num_days_min_dist_less_0.2m"
$ws.Range("M1").Value = "first_date_over_25"
$ws.Range("N1").Value = "consecutive_days_25"
$ws.Range("O1").Value = "first_date_over_100"
$ws.Range("P1").Value = "consecutive_days_100"
$ws.Range("Q1").Value = "n_poor_tracking_days"
$ws.Range("R1").Value = "total_missing_days"

# --- Step 2: Prepare date-format cells (M and O columns) by copying format from an existing date cell (E2) ---
$ws.Range("E2").Copy() | Out-Null
$ws.Range("M2:M19").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("O2:O19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 3: Clear any residual values pasted along with format assumptions; set blank-but-styled cells ---

# --- Step 4: Set data cell values row by row ---

# Row 2
$ws.Range("J2").ClearContents() | Out-Null
$ws.Range("K2").ClearContents() | Out-Null
$ws.Range("L2").ClearContents() | Out-Null
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").ClearContents() | Out-Null
$ws.Range("O2").ClearContents() | Out-Null
$ws.Range("P2").ClearContents() | Out-Null
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 17

# Row 3
$ws.Range("J3").ClearContents() | Out-Null
$ws.Range("K3").ClearContents() | Out-Null
$ws.Range("L3").ClearContents() | Out-Null
$ws.Range("M3").ClearContents() | Out-Null
$ws.Range("N3").ClearContents() | Out-Null
$ws.Range("O3").ClearContents() | Out-Null
$ws.Range("P3").ClearContents() | Out-Null
$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 20

# Row 4
$ws.Range("J4").ClearContents() | Out-Null
$ws.Range("K4").ClearContents() | Out-Null
$ws.Range("L4").ClearContents() | Out-Null
$ws.Range("M4").ClearContents() | Out-Null
$ws.Range("N4").ClearContents() | Out-Null
$ws.Range("O4").ClearContents() | Out-Null
$ws.Range("P4").ClearContents() | Out-Null
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = 0

# Row 5
$ws.Range("J5").ClearContents() | Out-Null
$ws.Range("K5").ClearContents() | Out-Null
$ws.Range("L5").ClearContents() | Out-Null
$ws.Range("M5").ClearContents() | Out-Null
$ws.Range("N5").ClearContents() | Out-Null
$ws.Range("O5").ClearContents() | Out-Null
$ws.Range("P5").ClearContents() | Out-Null
$ws.Range("Q5").Value = 11
$ws.Range("R5").Value = 16

# Row 6
$ws.Range("J6").ClearContents() | Out-Null
$ws.Range("K6").ClearContents() | Out-Null
$ws.Range("L6").ClearContents() | Out-Null
$ws.Range("M6").ClearContents() | Out-Null
$ws.Range("N6").ClearContents() | Out-Null
$ws.Range("O6").ClearContents() | Out-Null
$ws.Range("P6").ClearContents() | Out-Null
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 21

# Row 7
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1
$ws.Range("M7").ClearContents() | Out-Null
$ws.Range("N7").Value = 0
$ws.Range("O7").ClearContents() | Out-Null
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 17
$ws.Range("R7").Value = 0

# Row 8
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1
$ws.Range("M8").ClearContents() | Out-Null
$ws.Range("N8").Value = 0
$ws.Range("O8").ClearContents() | Out-Null
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0

# Row 9
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2
$ws.Range("M9").ClearContents() | Out-Null
$ws.Range("N9").Value = 0
$ws.Range("O9").ClearContents() | Out-Null
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 8
$ws.Range("R9").Value = 0

# Row 10
$ws.Range("J10").ClearContents() | Out-Null
$ws.Range("K10").ClearContents() | Out-Null
$ws.Range("L10").ClearContents() | Out-Null
$ws.Range("M10").ClearContents() | Out-Null
$ws.Range("N10").ClearContents() | Out-Null
$ws.Range("O10").ClearContents() | Out-Null
$ws.Range("P10").ClearContents() | Out-Null
$ws.Range("Q10").Value = 14
$ws.Range("R10").Value = 0

# Row 11
$ws.Range("J11").ClearContents() | Out-Null
$ws.Range("K11").ClearContents() | Out-Null
$ws.Range("L11").ClearContents() | Out-Null
$ws.Range("M11").ClearContents() | Out-Null
$ws.Range("N11").ClearContents() | Out-Null
$ws.Range("O11").ClearContents() | Out-Null
$ws.Range("P11").ClearContents() | Out-Null
$ws.Range("Q11").Value = 9
$ws.Range("R11").Value = 0

# Row 12
$ws.Range("J12").ClearContents() | Out-Null
$ws.Range("K12").ClearContents() | Out-Null
$ws.Range("L12").ClearContents() | Out-Null
$ws.Range("M12").ClearContents() | Out-Null
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("O12").ClearContents() | Out-Null
$ws.Range("P12").ClearContents() | Out-Null
$ws.Range("Q12").Value = 14
$ws.Range("R12").Value = 4

# Row 13
$ws.Range("J13").ClearContents() | Out-Null
$ws.Range("K13").ClearContents() | Out-Null
$ws.Range("L13").ClearContents() | Out-Null
$ws.Range("M13").ClearContents() | Out-Null
$ws.Range("N13").ClearContents() | Out-Null
$ws.Range("O13").ClearContents() | Out-Null
$ws.Range("P13").ClearContents() | Out-Null
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0

# Row 14
$ws.Range("J14").Value = 13
$ws.Range("K14").Value = 13
$ws.Range("L14").Value = 12
$ws.Range("M14").Value = 44836
$ws.Range("N14").Value = 13
$ws.Range("O14").Value = 44836
$ws.Range("P14").Value = 13
$ws.Range("Q14").Value = 27
$ws.Range("R14").Value = 0

# Row 15
$ws.Range("J15").ClearContents() | Out-Null
$ws.Range("K15").ClearContents() | Out-Null
$ws.Range("L15").ClearContents() | Out-Null
$ws.Range("M15").ClearContents() | Out-Null
$ws.Range("N15").ClearContents() | Out-Null
$ws.Range("O15").ClearContents() | Out-Null
$ws.Range("P15").ClearContents() | Out-Null
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 0

# Row 16
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 44841
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 44841
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = 13

# Row 17
$ws.Range("J17").ClearContents() | Out-Null
$ws.Range("K17").ClearContents() | Out-Null
$ws.Range("L17").ClearContents() | Out-Null
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").ClearContents() | Out-Null
$ws.Range("O17").ClearContents() | Out-Null
$ws.Range("P17").ClearContents() | Out-Null
$ws.Range("Q17").Value = 8
$ws.Range("R17").Value = 21

# Row 18
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 1
$ws.Range("M18").ClearContents() | Out-Null
$ws.Range("N18").Value = 0
$ws.Range("O18").ClearContents() | Out-Null
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 0

# Row 19
$ws.Range("J19").ClearContents() | Out-Null
$ws.Range("K19").ClearContents() | Out-Null
$ws.Range("L19").ClearContents() | Out-Null
$ws.Range("M19").ClearContents() | Out-Null
$ws.Range("N19").ClearContents() | Out-Null
$ws.Range("O19").ClearContents() | Out-Null
$ws.Range("P19").ClearContents() | Out-Null
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 0
